$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B, C, D columns with recalculated results for rows 2-10
$ws.Range("B2").Value = 0.9722
$ws.Range("C2").Value = 9.760400000000001
$ws.Range("D2").Value = 0.6955

$ws.Range("B3").Value = 0.7103
$ws.Range("C3").Value = 3.0618
$ws.Range("D3").Value = 0.6706

$ws.Range("B4").Value = 0.9172
$ws.Range("C4").Value = 10.969
$ws.Range("D4").Value = 0.646

$ws.Range("B5").Value = 0.5779
$ws.Range("C5").Value = 3.1462
$ws.Range("D5").Value = 0.6249

$ws.Range("B6").Value = 1.0382
$ws.Range("C6").Value = 5.8306
$ws.Range("D6").Value = 0.8408

$ws.Range("B7").Value = 0.7696
$ws.Range("C7").Value = 1.114
$ws.Range("D7").Value = 0.6751

$ws.Range("B8").Value = 1.0604
$ws.Range("C8").Value = 3.5094
$ws.Range("D8").Value = 0.6868

$ws.Range("B9").Value = 0.7222
$ws.Range("C9").Value = 2.8994
$ws.Range("D9").Value = 0.6916

$ws.Range("B10").Value = 0.5547
$ws.Range("C10").Value = 0.5479000000000001
$ws.Range("D10").Value = 0.6304

# LSPMW row's best scenario changed from ARMA to ARIMA
$ws.Range("E10").Value = "ARIMA"
